$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

$ws.Cells.Item(1, 10).Value = 62.2268602848053
$ws.Cells.Item(2, 2).Value = 1866
$ws.Cells.Item(2, 4).Value = 1863
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 99.94635193133047
$ws.Cells.Item(2, 8).Value = 99.89276139410188
$ws.Cells.Item(2, 9).Value = 0.00160857908847185
$ws.Cells.Item(2, 10).Value = 87.54390788078308
$ws.Cells.Item(3, 2).Value = 2083
$ws.Cells.Item(3, 4).Value = 2081
$ws.Cells.Item(3, 6).Value = 2
$ws.Cells.Item(3, 7).Value = 99.90398463754201
$ws.Cells.Item(3, 8).Value = 99.95196926032661
$ws.Cells.Item(3, 9).Value = 0.001439539347408829
$ws.Cells.Item(3, 10).Value = 87.76917457580566
$ws.Cells.Item(4, 2).Value = 2589
$ws.Cells.Item(4, 4).Value = 2562
$ws.Cells.Item(4, 6).Value = 9
$ws.Cells.Item(4, 7).Value = 99.64994165694283
$ws.Cells.Item(4, 8).Value = 98.99536321483771
$ws.Cells.Item(4, 9).Value = 0.01360808709175739
$ws.Cells.Item(4, 10).Value = 69.83948183059692
$ws.Cells.Item(5, 2).Value = 2026
$ws.Cells.Item(5, 4).Value = 2024
$ws.Cells.Item(5, 5).Value = 1
$ws.Cells.Item(5, 6).Value = 2
$ws.Cells.Item(5, 7).Value = 99.90128331688055
$ws.Cells.Item(5, 8).Value = 99.95061728395062
$ws.Cells.Item(5, 9).Value = 0.001480019733596448
$ws.Cells.Item(5, 10).Value = 87.93068289756775
$ws.Cells.Item(6, 10).Value = 74.98183512687683
$ws.Cells.Item(7, 2).Value = 2530
$ws.Cells.Item(7, 4).Value = 2529
$ws.Cells.Item(7, 6).Value = 2
$ws.Cells.Item(7, 7).Value = 99.92097984986171
$ws.Cells.Item(7, 9).Value = 0.0007898894154818325
$ws.Cells.Item(7, 10).Value = 103.6006164550781
$ws.Cells.Item(8, 10).Value = 84.08601069450378
$ws.Cells.Item(9, 10).Value = 107.9590029716492
$ws.Cells.Item(10, 2).Value = 1795
$ws.Cells.Item(10, 5).Value = 0
$ws.Cells.Item(10, 8).Value = 100
$ws.Cells.Item(10, 9).Value = 0
$ws.Cells.Item(10, 10).Value = 85.13120865821838
$ws.Cells.Item(11, 10).Value = 65.15259194374084
$ws.Cells.Item(12, 10).Value = 83.35275340080261
$ws.Cells.Item(13, 2).Value = 2389
$ws.Cells.Item(13, 4).Value = 2387
$ws.Cells.Item(13, 5).Value = 1
$ws.Cells.Item(13, 6).Value = 24
$ws.Cells.Item(13, 7).Value = 99.00456242223144
$ws.Cells.Item(13, 8).Value = 99.95812395309883
$ws.Cells.Item(13, 9).Value = 0.01036484245439469
$ws.Cells.Item(13, 10).Value = 89.46774291992188
$ws.Cells.Item(14, 10).Value = 77.64330434799194
$ws.Cells.Item(15, 2).Value = 2280
$ws.Cells.Item(15, 5).Value = 2
$ws.Cells.Item(15, 8).Value = 99.91224221149628
$ws.Cells.Item(15, 9).Value = 0.000877963125548727
$ws.Cells.Item(15, 10).Value = 93.17987847328186
$ws.Cells.Item(16, 2).Value = 2021
$ws.Cells.Item(16, 5).Value = 34
$ws.Cells.Item(16, 8).Value = 98.31683168316832
$ws.Cells.Item(16, 9).Value = 0.0171112229491696
$ws.Cells.Item(16, 10).Value = 114.8682944774628
$ws.Cells.Item(17, 10).Value = 119.6526062488556
$ws.Cells.Item(18, 10).Value = 151.6958434581757
$ws.Cells.Item(19, 10).Value = 112.1011893749237
$ws.Cells.Item(20, 10).Value = 70.75770592689514
$ws.Cells.Item(21, 10).Value = 118.4828763008118
$ws.Cells.Item(22, 10).Value = 141.2787873744965
$ws.Cells.Item(23, 2).Value = 2123
$ws.Cells.Item(23, 4).Value = 2121
$ws.Cells.Item(23, 6).Value = 14
$ws.Cells.Item(23, 7).Value = 99.34426229508196
$ws.Cells.Item(23, 9).Value = 0.006554307116104869
$ws.Cells.Item(23, 10).Value = 81.79234981536865
$ws.Cells.Item(24, 2).Value = 2964
$ws.Cells.Item(24, 4).Value = 2935
$ws.Cells.Item(24, 5).Value = 28
$ws.Cells.Item(24, 6).Value = 44
$ws.Cells.Item(24, 7).Value = 98.52299429338704
$ws.Cells.Item(24, 8).Value = 99.05501181235235
$ws.Cells.Item(24, 9).Value = 0.02416107382550335
$ws.Cells.Item(24, 10).Value = 162.7516641616821
$ws.Cells.Item(25, 2).Value = 2649
$ws.Cells.Item(25, 4).Value = 2648
$ws.Cells.Item(25, 6).Value = 7
$ws.Cells.Item(25, 7).Value = 99.73634651600753
$ws.Cells.Item(25, 9).Value = 0.002635542168674699
$ws.Cells.Item(25, 10).Value = 168.0317049026489
$ws.Cells.Item(26, 10).Value = 141.738062620163
$ws.Cells.Item(27, 2).Value = 2944
$ws.Cells.Item(27, 4).Value = 2940
$ws.Cells.Item(27, 5).Value = 3
$ws.Cells.Item(27, 6).Value = 14
$ws.Cells.Item(27, 7).Value = 99.52606635071091
$ws.Cells.Item(27, 8).Value = 99.8980632008155
$ws.Cells.Item(27, 9).Value = 0.005752961082910322
$ws.Cells.Item(27, 10).Value = 151.8597741127014
$ws.Cells.Item(28, 10).Value = 157.9463520050049
$ws.Cells.Item(29, 2).Value = 2622
$ws.Cells.Item(29, 4).Value = 2619
$ws.Cells.Item(29, 5).Value = 2
$ws.Cells.Item(29, 6).Value = 30
$ws.Cells.Item(29, 7).Value = 98.86749716874291
$ws.Cells.Item(29, 8).Value = 99.92369324685235
$ws.Cells.Item(29, 9).Value = 0.01207547169811321
$ws.Cells.Item(29, 10).Value = 157.0292055606842
$ws.Cells.Item(30, 10).Value = 165.1074182987213
$ws.Cells.Item(31, 2).Value = 3248
$ws.Cells.Item(31, 4).Value = 3247
$ws.Cells.Item(31, 5).Value = 0
$ws.Cells.Item(31, 6).Value = 3
$ws.Cells.Item(31, 7).Value = 99.9076923076923
$ws.Cells.Item(31, 8).Value = 100
$ws.Cells.Item(31, 9).Value = 0.0009227929867733005
$ws.Cells.Item(31, 10).Value = 222.9341397285461
$ws.Cells.Item(32, 2).Value = 2259
$ws.Cells.Item(32, 5).Value = 0
$ws.Cells.Item(32, 8).Value = 100
$ws.Cells.Item(32, 9).Value = 0.001326259946949602
$ws.Cells.Item(32, 10).Value = 173.94384765625
$ws.Cells.Item(33, 10).Value = 193.7815217971802
$ws.Cells.Item(34, 2).Value = 2155
$ws.Cells.Item(34, 4).Value = 2152
$ws.Cells.Item(34, 6).Value = 1
$ws.Cells.Item(34, 7).Value = 99.95355318160706
$ws.Cells.Item(34, 8).Value = 99.90714948932219
$ws.Cells.Item(34, 9).Value = 0.001392757660167131
$ws.Cells.Item(34, 10).Value = 180.2487514019012
$ws.Cells.Item(35, 10).Value = 248.1632480621338
$ws.Cells.Item(36, 10).Value = 170.9635078907013
$ws.Cells.Item(37, 10).Value = 173.4883069992065
$ws.Cells.Item(38, 2).Value = 2600
$ws.Cells.Item(38, 4).Value = 2598
$ws.Cells.Item(38, 6).Value = 6
$ws.Cells.Item(38, 7).Value = 99.76958525345623
$ws.Cells.Item(38, 9).Value = 0.002303262955854127
$ws.Cells.Item(38, 10).Value = 132.4153089523315
$ws.Cells.Item(39, 2).Value = 2058
$ws.Cells.Item(39, 4).Value = 2049
$ws.Cells.Item(39, 5).Value = 8
$ws.Cells.Item(39, 6).Value = 3
$ws.Cells.Item(39, 7).Value = 99.85380116959064
$ws.Cells.Item(39, 8).Value = 99.61108410306271
$ws.Cells.Item(39, 10).Value = 214.5361526012421
$ws.Cells.Item(40, 10).Value = 248.2104690074921
$ws.Cells.Item(41, 10).Value = 151.5301859378815
$ws.Cells.Item(42, 10).Value = 126.1397602558136
$ws.Cells.Item(43, 10).Value = 266.1511223316193
$ws.Cells.Item(44, 2).Value = 2752
$ws.Cells.Item(44, 4).Value = 2751
$ws.Cells.Item(44, 6).Value = 1
$ws.Cells.Item(44, 7).Value = 99.96366279069767
$ws.Cells.Item(44, 9).Value = 0.0003632401017072285
$ws.Cells.Item(44, 10).Value = 249.2073049545288
